$d = $word.ActiveDocument

# --- 1. "JOHP;" paragraph: drop the proofErr wrappers and merge the two
#        runs ("JOHP;" + " ") into a single run "JOHP; ".
$pJohp = $d.Paragraphs.Item(3)
$pJohp.Range.InsertParagraphBefore()
$d.Paragraphs.Item(3).Range.Text = "JOHP; "
$d.Paragraphs.Item(4).Range.Delete()

# --- 2. "work and stress;" paragraph: drop the proofErr wrappers and merge
#        the two runs ("work and " + "stress;") into a single run.
$pStress = $d.Paragraphs.Item(4)
$pStress.Range.InsertParagraphBefore()
$d.Paragraphs.Item(4).Range.Text = "work and stress;"
$d.Paragraphs.Item(5).Range.Delete()

# --- 3. Insert a new paragraph "International Journal of Stress Management"
#        right after "European Journal of Work and Organizational Psychology".
$pEuro = $d.Paragraphs.Item(6)
$rEuro = $pEuro.Range
$rEuro.Collapse(0)
$rEuro.InsertParagraphAfter()
$d.Paragraphs.Item(7).Range.Text = "International Journal of Stress Management"

# --- 4. The old "International Journal of Stress Management" paragraph
#        (now the last paragraph, under "Potential Options") becomes the
#        "**Could consider..." note, followed by a long list of new outlets.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$pLast.Range.Text = "**Could consider HR-focused journals too…"

$newOutlets = @(
    "Scandinavian Journal of Work, Environment & Health",
    "Personality and Individual Differences",
    "Social Science Research",
    "Industrial and Organizational Psychology",
    "Journal of Individual Differences",
    "The Psychologist-Manager Journal,",
    "Journal of Vocational Education & Training",
    "Journal of Managerial Psychology",
    "Employee Relations",
    "Applied Psychology",
    "Journal of Occupational and Organizational Psychology",
    "Journal of Management",
    "Work and Occupations",
    "Journal of Career Development",
    "Journal of Management Education",
    "Journal of Career Assessment",
    "Management Science",
    "Work, 41(Suppl. 1)",
    "Occupational Outlook Quarterly"
)

$cur = $d.Paragraphs.Item($d.Paragraphs.Count)
foreach ($outlet in $newOutlets) {
    $cur.Range.InsertParagraphAfter()
    $cur = $d.Paragraphs.Item($d.Paragraphs.Count)
    $cur.Range.Text = $outlet
}

# "Journal of Business and Psychology" is split across two runs in the
# source document, so build it as two separate InsertAfter calls and then
# stitch the two paragraphs together by deleting the paragraph mark between
# them (rather than letting Range.Text silently merge the text into a
# single run).
$cur.Range.InsertParagraphAfter()
$pBiz1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$pBiz1.Range.Text = "Journal of Business and"

$pBiz1.Range.InsertParagraphAfter()
$pBiz2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$pBiz2.Range.Text = " Psychology"

$markPos = $pBiz1.Range.End - 1
$d.Range($markPos, $markPos + 1).Delete()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count)

$moreOutlets = @(
    "Journal of Human Resources Education",
    "Journal of Career Development",
    "New Zealand Journal of Psychology"
)

foreach ($outlet in $moreOutlets) {
    $cur.Range.InsertParagraphAfter()
    $cur = $d.Paragraphs.Item($d.Paragraphs.Count)
    $cur.Range.Text = $outlet
}
